# Renamed few transcripts in the "Speaker" column (column D) of the DataSheet.
# "RBD" -> "T" and "Student" -> "S" for every matching cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)  # Column D = Speaker
    $val = $cell.Value2
    if ($val -eq "RBD") {
        $cell.Value2 = "T"
    } elseif ($val -eq "Student") {
        $cell.Value2 = "S"
    }
}
